# Update "本次参展人数" (F column) figures across the four sheets of the
# workbook, matching the regenerated output published to gh-pages.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 -----------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$updates = @{
    2  = 7709
    3  = 100
    4  = 78
    5  = 7555
    7  = 591
    8  = 627
    9  = 446
    11 = 437
    12 = 772
    13 = 35
    15 = 299
    17 = 259
    18 = 136
    19 = 391
    20 = 145
    22 = 76
    23 = 602
    24 = 2196
    25 = 723
    26 = 49
    27 = 50
    29 = 607
    30 = 51
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# --- Sheet: 演出 -----------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$updates = @{
    2  = 286
    4  = 321
    5  = 315
    10 = 1
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# --- Sheet: 本地生活 ---------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$updates = @{
    2 = 443
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# --- Sheet: 全部类型 ---------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$updates = @{
    2  = 443
    3  = 7709
    4  = 100
    5  = 78
    6  = 286
    7  = 7556
    9  = 591
    10 = 627
    11 = 446
    14 = 437
    15 = 321
    16 = 315
    18 = 772
    19 = 35
    21 = 299
    26 = 259
    27 = 136
    28 = 391
    29 = 145
    31 = 76
    32 = 602
    33 = 2196
    34 = 723
    35 = 49
    36 = 50
    38 = 1
    39 = 607
    40 = 51
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
